$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 382, shifting existing rows 382..470 down to 383..471
$ws.Rows.Item(382).Insert()

# Populate the newly inserted row 382 with the new data record
$ws.Cells.Item(382, 1).Value = 2
$ws.Cells.Item(382, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(382, 3).Value = "Coquimbo"
$ws.Cells.Item(382, 4).Value = 45204
$ws.Cells.Item(382, 5).Value = 4
$ws.Cells.Item(382, 6).Value = 100112021
$ws.Cells.Item(382, 7).Value = "Ají"
$ws.Cells.Item(382, 8).Value = "Americana (o)"
$ws.Cells.Item(382, 9).Value = "Primera"
$ws.Cells.Item(382, 10).Value = 200
$ws.Cells.Item(382, 11).Value = 60000
$ws.Cells.Item(382, 12).Value = 65000
$ws.Cells.Item(382, 13).Value = 62500
$ws.Cells.Item(382, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(382, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(382, 16).Value = 2500
$ws.Cells.Item(382, 17).Value = 25
$ws.Cells.Item(382, 18).Value = "Hortaliza"
